# "Generate Report for Handback" -- the handoff/handback tracker report gets
# refreshed: the two in-flight files (f13d269d... and ffffb80e2993...) have now
# been handed back and are in sync with en-US, so:
#   * the Status column flips from "Ready for handoff" to
#     "Handed back: in sync with en-US" on every locale sheet
#   * each locale sheet grows two columns ("Latest Target File" / E, and
#     "Latest Handback File" / F) recording the handback artifacts, with
#     their own hyperlinks
#   * "Latest Handback DateTime" (column G) is stamped with the handback time
#     for each locale (zh-cn handed back a few seconds before de-de)

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$mdDisplay = "f13d269d-2389-4baf-9322-ab170051d945.md"
$mdUrl     = "https://github.com/OpenLocalizationTest/oltest/blob/ab7c25a70fe6a49b0d7abce5bb598014b6520bd0/e2e/f13d269d-2389-4baf-9322-ab170051d945.md"

function Update-LocaleSheet {
    param($SheetName, $XlfDisplay, $XlfUrl, $HandbackTime)

    $ws = $wb.Worksheets.Item($SheetName)

    foreach ($row in 2, 3) {
        # Status: handoff is complete, now in sync with en-US.
        $ws.Range("B$row").Value = $statusText

        # E: Latest Target File -- the handed-back .md, linked back to source.
        $ws.Hyperlinks.Add($ws.Range("E$row"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdDisplay)

        # F: Latest Handback File -- the handed-back xlf for this locale.
        $ws.Hyperlinks.Add($ws.Range("F$row"), $XlfUrl, [Type]::Missing, [Type]::Missing, $XlfDisplay)

        # G: Latest Handback DateTime -- now stamped instead of the epoch sentinel.
        $ws.Range("G$row").Value = $HandbackTime
    }
}

Update-LocaleSheet "zh-cn" `
    "f13d269d-2389-4baf-9322-ab170051d945.b476a835270279f96d7c74b645f3371bdbfad9e9.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ec4d7c626e85b0a8e1d3603047f57a38cc8b313b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/f13d269d-2389-4baf-9322-ab170051d945.b476a835270279f96d7c74b645f3371bdbfad9e9.zh-cn.xlf" `
    "2016-02-22 17:57:26"

Update-LocaleSheet "de-de" `
    "f13d269d-2389-4baf-9322-ab170051d945.b476a835270279f96d7c74b645f3371bdbfad9e9.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e50deee38aebe23fc56cea68436aa17ad67e82f0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/f13d269d-2389-4baf-9322-ab170051d945.b476a835270279f96d7c74b645f3371bdbfad9e9.de-de.xlf" `
    "2016-02-22 17:57:45"

# The "Overview" sheet's zh-cn/de-de status columns share the very same
# "Ready for handoff" string as the per-locale Status column, so they flip
# to the new status text too.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusText
$wsOverview.Range("C2").Value = $statusText
$wsOverview.Range("B3").Value = $statusText
$wsOverview.Range("C3").Value = $statusText
